$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Periodo Mora" column (E16:E22) is re-ordered from descending (1912..1906)
# to ascending (1906..1912) as part of the EC database refresh / "parte 1"
# of new estado de cuenta rows.
$ws.Range("E16").Value = "1906"
$ws.Range("E17").Value = "1907"
$ws.Range("E18").Value = "1908"
$ws.Range("E19").Value = "1909"
$ws.Range("E20").Value = "1910"
$ws.Range("E21").Value = "1911"
$ws.Range("E22").Value = "1912"

# "Valor Mora" column (F16:F22) keeps the same set of values but follows the
# row re-ordering: the odd 20267 value moves from row 16 to row 22.
$ws.Range("F16").Value = 38000
$ws.Range("F17").Value = 38000
$ws.Range("F18").Value = 38000
$ws.Range("F19").Value = 38000
$ws.Range("F20").Value = 38000
$ws.Range("F21").Value = 38000
$ws.Range("F22").Value = 20267
